$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text representation
# (values like "258.40" or "2.27" would otherwise be auto-converted to numbers)
$numRange = $ws.Range("D2:E51")
$numRange.NumberFormat = "@"

# --- Rows with simple price/volume refreshes ---
$ws.Range("D2").Value = '43.666.25'
$ws.Range("E2").Value = '  +1.73%  '
$ws.Range("D3").Value = '2.204.10'
$ws.Range("E3").Value = '  -0.58%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '258.40'
$ws.Range("E5").Value = '  +1.64%  '
$ws.Range("D6").Value = '85.08'
$ws.Range("E6").Value = '  +11.90%  '
$ws.Range("D7").Value = '0.618'
$ws.Range("E7").Value = '  +0.59%  '
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("D9").Value = '0.594'
$ws.Range("E9").Value = '  +0.46%  '
$ws.Range("D10").Value = '44.82'
$ws.Range("E10").Value = '  +8.38%  '
$ws.Range("D11").Value = '0.0915'
$ws.Range("E11").Value = '  +0.25%  '
$ws.Range("D12").Value = '7.37'
$ws.Range("E12").Value = '  +6.94%  '
$ws.Range("E13").Value = '  +1.67%  '
$ws.Range("D14").Value = '2.532.91'
$ws.Range("E14").Value = '  -0.77%  '
$ws.Range("E15").Value = '  -0.40%  '
$ws.Range("D16").Value = '2.209.72'
$ws.Range("E16").Value = '  -0.31%  '
$ws.Range("D17").Value = '0.782'
$ws.Range("E17").Value = '  -0.50%  '
$ws.Range("D18").Value = '43.632.21'
$ws.Range("E18").Value = '  +1.80%  '
$ws.Range("E19").Value = '  +0.42%  '
$ws.Range("D20").Value = '69.79'
$ws.Range("E20").Value = '  -1.94%  '
$ws.Range("D21").Value = '5.91'
$ws.Range("E21").Value = '  -0.56%  '
$ws.Range("E22").Value = '  +6.29%  '
$ws.Range("D23").Value = '230.92'
$ws.Range("E23").Value = '  +0.68%  '
$ws.Range("D24").Value = '8.93'
$ws.Range("E24").Value = '  -4.25%  '
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("D26").Value = '3.56'
$ws.Range("E26").Value = '  +5.77%  '
$ws.Range("D27").Value = '10.62'
$ws.Range("E27").Value = '  +0.14%  '
$ws.Range("D30").Value = '2.24'
$ws.Range("E30").Value = '  +2.42%  '
$ws.Range("D31").Value = '173.73'
$ws.Range("D32").Value = '20.39'
$ws.Range("E32").Value = '  +0.75%  '
$ws.Range("D33").Value = '0.0866'
$ws.Range("E33").Value = '  +1.82%  '
$ws.Range("D34").Value = '5.33'
$ws.Range("E34").Value = '  +2.49%  '
$ws.Range("D35").Value = '0.123'
$ws.Range("E35").Value = '  +1.32%  '
$ws.Range("E36").Value = '  +0.59%  '
$ws.Range("D39").Value = '12.48'
$ws.Range("E39").Value = '  +0.45%  '
$ws.Range("D40").Value = '2.85'
$ws.Range("E40").Value = '  +3.41%  '
$ws.Range("E41").Value = '  -0.51%  '
$ws.Range("E42").Value = '  +5.35%  '
$ws.Range("D43").Value = '5.46'
$ws.Range("E43").Value = '  +3.67%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D47").Value = '0.0977'
$ws.Range("E47").Value = '  -0.36%  '
$ws.Range("E48").Value = '  +4.11%  '
$ws.Range("D49").Value = '1.11'
$ws.Range("E49").Value = '  +0.19%  '
$ws.Range("D50").Value = '0.438'
$ws.Range("E50").Value = '  -3.54%  '
$ws.Range("D51").Value = '1.47'
$ws.Range("E51").Value = '  +1.75%  '

# --- Rows whose coin data was reordered (content swapped between two rows) ---
$ws.Range("B28").Value = 'InjectiveProtocol'
$ws.Range("C28").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D28").Value = '39.38'
$ws.Range("E28").Value = '  +1.77%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '2.27'
$ws.Range("E29").Value = '  +2.75%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.0359'
$ws.Range("E37").Value = '  +1.78%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '4.46'
$ws.Range("E38").Value = '  +3.35%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '100.35'
$ws.Range("E45").Value = '  -1.35%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '8.34'
$ws.Range("E46").Value = '  -0.04%  '

# Restore normal (default) cell style now that the text values are set,
# matching the original workbook formatting (no explicit number format).
$numRange.Style = "Normal"
